# SEMANA1.xlsx edit: replace the three text "fecha" columns (stored as
# shared-string dd/mmm/yyyy labels with a text number format) with real
# Excel date serials formatted as dd-mm-yy, and move the selection/scroll
# position to where the author left off (around row 78).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dateFormat = "dd-mm-yy;@"

# Block 1: rows 2-39 -> 19/abr/2021 (serial 44305)
$r1 = $ws.Range("A2:A39")
$r1.NumberFormat = $dateFormat
$r1.Value2 = 44305

# Block 2: rows 40-77 -> 18/abr/2021 (serial 44304)
$r2 = $ws.Range("A40:A77")
$r2.NumberFormat = $dateFormat
$r2.Value2 = 44304

# Block 3: rows 78-115 -> 17/abr/2021 (serial 44303)
$r3 = $ws.Range("A78:A115")
$r3.NumberFormat = $dateFormat
$r3.Value2 = 44303

# Restore the author's last selection / scroll position (row 78 area).
$ws.Range("A78").Select()
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
